# docs+wbs: close partials with CI evidence and category-gated done
#
# 1. On the WBS sheet, rows 2-49 (all rows whose Execution Status was
#    "Done") are downgraded to "Partial" and their "Completed On" (K)
#    date is cleared, since completion is now gated on the new
#    category columns.
# 2. Five new governance-category columns are appended as headers in
#    row 1: Schema, Validation, Permissions/Isolation, Workflow,
#    Evidence (columns L-P), expanding the used range to A1:P137.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WBS")

# --- New category-gate header columns (L1:P1) ---------------------------
$ws.Range("L1").Value = "Schema"
$ws.Range("M1").Value = "Validation"
$ws.Range("N1").Value = "Permissions/Isolation"
$ws.Range("O1").Value = "Workflow"
$ws.Range("P1").Value = "Evidence"

# --- Downgrade previously "Done" rows (2-49) to "Partial" and clear the
#     Completed On (K) date, since completion now requires the new
#     category gates above to be satisfied. ----------------------------
for ($r = 2; $r -le 49; $r++) {
    $ws.Cells.Item($r, 8).Value = "Partial"
    $ws.Cells.Item($r, 11).ClearContents()
}
